$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.985.45"
$ws.Range("E2").Value = "  -0.67%  "

# Row 3
$ws.Range("D3").Value = "2.514.33"
$ws.Range("E3").Value = "  -0.55%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.88%  "

# Row 7
$ws.Range("E7").Value = "  -0.14%  "

# Row 9
$ws.Range("D9").Value = "2.512.97"
$ws.Range("E9").Value = "  -0.77%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.28%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.66%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.38%  "

# Row 14
$ws.Range("D14").Value = "2.935.60"
$ws.Range("E14").Value = "  -1.31%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "

# Row 16
$ws.Range("D16").Value = "58.725.92"
$ws.Range("E16").Value = "  -1.00%  "

# Row 18
$ws.Range("D18").Value = "2.504.19"
$ws.Range("E18").Value = "  -2.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "

# Row 20
$ws.Range("E20").Value = "  +0.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "

# Row 22
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("E23").Value = "  +0.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.41%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.418"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "

# Row 26
$ws.Range("E26").Value = "  -0.31%  "

# Row 27
$ws.Range("E27").Value = "  -0.16%  "

# Row 28
$ws.Range("E28").Value = "  -3.65%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "

# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  -0.08%  "

# Row 31
$ws.Range("E31").Value = "  -2.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.67%  "

# Row 33
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("E34").Value = "  +0.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.85%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.79%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "

# Row 41
$ws.Range("E41").Value = "  -1.23%  "

# Row 42
$ws.Range("E42").Value = "  -0.72%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "

# Row 45
$ws.Range("E45").Value = "  -0.28%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.32%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0926"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.80%  "

# Row 50
$ws.Range("E50").Value = "  -1.05%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.72%  "
